# Update TPM-derived NATMI metrics in the active worksheet.
# The diff updates ligand/receptor average & total expression values
# (columns G,H,I,J for ligand; M,N,O,P for receptor) and the derived
# edge weights/specificities (columns Q,R,S,T) for rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ G = 1.107578333333333;  H = 3.322735;   I = 0.004477262396335327; J = 0.004477262396335327;
            M = 3.456265333333333;  N = 10.368796;  O = 0.009841535807677501; P = 0.0098415358076775;
            Q = 3.828084597451111;  R = 34.45276137706; S = 0.0000440631381939021; T = 0.0000440631381939021 }
    3  = @{ G = 1.107578333333333;  H = 3.322735;   I = 0.004477262396335327; J = 0.004477262396335327;
            O = 0.8587907398420774; P = 0.8587907398420773;
            Q = 334.0457899933122;  R = 3006.41210993981; S = 0.003845031485815928; T = 0.003845031485815928 }
    4  = @{ G = 1.107578333333333;  H = 3.322735;   I = 0.004477262396335327; J = 0.004477262396335327;
            O = 0.1313677243502452; P = 0.1313677243502452;
            Q = 51.0984029337239;   R = 459.885626403515; S = 0.0005881677723254975; T = 0.0005881677723254974 }
    5  = @{ G = 212.7693433333334;  H = 638.30803;  I = 0.860096438625976;    J = 0.8600964386259761;
            M = 3.456265333333333;  N = 10.368796;  O = 0.009841535807677501; P = 0.0098415358076775;
            Q = 735.3873053590978;  R = 6618.48574823188; S = 0.008464669898793437; T = 0.008464669898793437 }
    6  = @{ G = 212.7693433333334;  H = 638.30803;  I = 0.860096438625976;    J = 0.8600964386259761;
            O = 0.8587907398420774; P = 0.8587907398420773;
            Q = 64171.26558104238;  R = 577541.3902293814; S = 0.7386428568631378; T = 0.7386428568631378 }
    7  = @{ G = 212.7693433333334;  H = 638.30803;  I = 0.860096438625976;    J = 0.8600964386259761;
            O = 0.1313677243502452; P = 0.1313677243502452;
            Q = 9816.166776096055;  R = 88345.50098486448; S = 0.1129889118640448; T = 0.1129889118640448 }
    8  = @{ I = 0.1354262989776887; J = 0.1354262989776887;
            M = 3.456265333333333;  N = 10.368796;  O = 0.009841535807677501; P = 0.0098415358076775;
            Q = 115.7902493341986;  R = 1042.112244007788; S = 0.001332802770690162; T = 0.001332802770690162 }
    9  = @{ I = 0.1354262989776887; J = 0.1354262989776887;
            O = 0.8587907398420774; P = 0.8587907398420773;
            S = 0.1163028514931237; T = 0.1163028514931237 }
    10 = @{ I = 0.1354262989776887; J = 0.1354262989776887;
            O = 0.1313677243502452; P = 0.1313677243502452;
            S = 0.0177906447138749; T = 0.0177906447138749 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $addr = "$colLetter$rowNum"
        $ws.Range($addr).Value = $cols[$colLetter]
    }
}
